# Updated cryptos list with latest fetch (Price + Volume(1h) columns).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "66.015.29"

$ws.Range("D3").Value = "3.445.25"
$ws.Range("E3").Value = "  -0.19%  "

$ws.Range("E4").Value = "  +0.00%  "

$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = "585.71"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.99%  "

$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value = "174.12"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -0.99%  "

$ws.Range("E7").Value = "  +0.02%  "

$c = $ws.Cells.Item(8, 4)
$c.NumberFormat = "@"
$c.Value = "0.603"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +0.57%  "

$ws.Range("D9").Value = "3.443.33"
$ws.Range("E9").Value = "  -0.18%  "

$ws.Range("E10").Value = "  -1.42%  "

$ws.Range("E11").Value = "  +1.39%  "

$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = "@"
$c.Value = "0.420"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +0.30%  "

$ws.Range("D13").Value = "4.042.56"
$ws.Range("E13").Value = "  -0.06%  "

$ws.Range("E14").Value = "  +1.81%  "

$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = "@"
$c.Value = "29.27"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -3.84%  "

$ws.Range("D16").Value = "65.978.66"
$ws.Range("E16").Value = "  -0.79%  "

$ws.Range("E17").Value = "  -0.09%  "

$ws.Range("D18").Value = "3.439.68"
$ws.Range("E18").Value = "  -0.30%  "

$ws.Range("E20").Value = "  -0.02%  "

$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = "@"
$c.Value = "370.74"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -1.39%  "

$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = "7.62"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -1.03%  "

$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = "@"
$c.Value = "72.56"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +2.01%  "

$ws.Range("E24").Value = "  -0.08%  "

$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = "@"
$c.Value = "0.534"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +1.31%  "

$ws.Range("E26").Value = "  +4.11%  "

$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = "@"
$c.Value = "9.73"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -0.44%  "

$ws.Range("E28").Value = "  +3.92%  "

$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -0.05%  "

$ws.Range("E30").Value = "  -0.64%  "

$ws.Range("E31").Value = "  +0.48%  "

$ws.Range("E32").Value = "  -1.40%  "

$ws.Range("E33").Value = "  -0.02%  "

$ws.Range("E34").Value = "  +0.04%  "

$ws.Range("E35").Value = "  -4.62%  "

$c = $ws.Cells.Item(36, 4)
$c.NumberFormat = "@"
$c.Value = "1.54"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +1.67%  "

$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = "@"
$c.Value = "161.78"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +1.62%  "

$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = "@"
$c.Value = "0.879"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +0.29%  "

$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = "@"
$c.Value = "28.37"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +3.66%  "

$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value = "1.79"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +0.33%  "

$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value = "2.61"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -0.80%  "

$ws.Range("D42").Value = "2.788.67"
$ws.Range("E42").Value = "  +3.95%  "

$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = "@"
$c.Value = "4.48"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +0.38%  "

$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = "@"
$c.Value = "6.46"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -0.63%  "

$ws.Range("E45").Value = "  -0.89%  "

$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = "@"
$c.Value = "25.24"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +0.22%  "

$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = "39.78"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -0.97%  "

$ws.Range("E48").Value = "  -0.93%  "

$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = "@"
$c.Value = "327.45"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +2.41%  "

$ws.Range("E50").Value = "  +0.39%  "

$ws.Range("E51").Value = "  +1.59%  "
